$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Citywide Totals")
Write-Output $ws.Name
$ws.Range("E2").Value = 27
